$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Formula = "'275.65"
$ws.Range("F2").Value = "13-12-2022"
$ws.Range("G2").Formula = "'0"

# Row 3
$ws.Range("D3").Formula = "'21.08"
$ws.Range("F3").Value = "13-12-2022"
$ws.Range("G3").Formula = "'0"

# Row 4
$ws.Range("D4").Formula = "'6.257"
$ws.Range("F4").Value = "13-12-2022"
$ws.Range("G4").Formula = "'0"

# Row 5
$ws.Range("D5").Formula = "'0.06210"
$ws.Range("F5").Value = "13-12-2022"
$ws.Range("G5").Formula = "'0"

# Row 6
$ws.Range("D6").Formula = "'3.557"
$ws.Range("F6").Value = "13-12-2022"
$ws.Range("G6").Formula = "'0"

# Row 7
$ws.Range("D7").Formula = "'6.585"
$ws.Range("F7").Value = "13-12-2022"
$ws.Range("G7").Formula = "'0"

# Row 8
$ws.Range("D8").Formula = "'1.482"
$ws.Range("F8").Value = "13-12-2022"
$ws.Range("G8").Formula = "'0"

# Row 9
$ws.Range("D9").Formula = "'0.8277"
$ws.Range("F9").Value = "13-12-2022"
$ws.Range("G9").Formula = "'0"

# Row 10
$ws.Range("D10").Formula = "'0.1666"
$ws.Range("F10").Value = "13-12-2022"
$ws.Range("G10").Formula = "'0"

# Row 11
$ws.Range("D11").Formula = "'0.08321"
$ws.Range("F11").Value = "13-12-2022"
$ws.Range("G11").Formula = "'0"

# Row 12
$ws.Range("D12").Formula = "'0.03524"
$ws.Range("F12").Value = "13-12-2022"
$ws.Range("G12").Formula = "'0"

# Row 13
$ws.Range("D13").Formula = "'0.03188"
$ws.Range("F13").Value = "13-12-2022"
$ws.Range("G13").Formula = "'0"

# Row 14
$ws.Range("D14").Formula = "'0.09171"
$ws.Range("F14").Value = "13-12-2022"
$ws.Range("G14").Formula = "'0"

# Row 15
$ws.Range("F15").Value = "13-12-2022"
$ws.Range("G15").Formula = "'0"

# Row 16
$ws.Range("D16").Formula = "'0.001640"
$ws.Range("F16").Value = "13-12-2022"
$ws.Range("G16").Formula = "'0"

# Row 17
$ws.Range("D17").Formula = "'0.04707"
$ws.Range("F17").Value = "13-12-2022"
$ws.Range("G17").Formula = "'0"

# Row 18
$ws.Range("D18").Formula = "'0.006327"
$ws.Range("F18").Value = "13-12-2022"
$ws.Range("G18").Formula = "'0"

# Row 19
$ws.Range("D19").Formula = "'0.006209"
$ws.Range("F19").Value = "13-12-2022"
$ws.Range("G19").Formula = "'0"

# Row 20
$ws.Range("F20").Value = "13-12-2022"
$ws.Range("G20").Formula = "'0"

# Row 21
$ws.Range("F21").Value = "13-12-2022"
$ws.Range("G21").Formula = "'0"

# Row 22
$ws.Range("D22").Formula = "'3.719"
$ws.Range("F22").Value = "13-12-2022"
$ws.Range("G22").Formula = "'0"

# Row 23
$ws.Range("D23").Formula = "'2.264"
$ws.Range("F23").Value = "13-12-2022"
$ws.Range("G23").Formula = "'0"

# Row 24
$ws.Range("D24").Formula = "'0.01396"
$ws.Range("F24").Value = "13-12-2022"
$ws.Range("G24").Formula = "'0"

# Row 25
$ws.Range("F25").Value = "13-12-2022"
$ws.Range("G25").Formula = "'0"

# Row 26
$ws.Range("F26").Value = "13-12-2022"
$ws.Range("G26").Formula = "'0"

# Row 27
$ws.Range("F27").Value = "13-12-2022"
$ws.Range("G27").Formula = "'0"

# Row 28
$ws.Range("D28").Formula = "'0.0002726"
$ws.Range("F28").Value = "13-12-2022"
$ws.Range("G28").Formula = "'0"

# Row 29
$ws.Range("F29").Value = "13-12-2022"
$ws.Range("G29").Formula = "'0"

# Row 30
$ws.Range("F30").Value = "13-12-2022"
$ws.Range("G30").Formula = "'0"

# Row 31
$ws.Range("F31").Value = "13-12-2022"
$ws.Range("G31").Formula = "'0"

# Row 32
$ws.Range("F32").Value = "13-12-2022"
$ws.Range("G32").Formula = "'0"

# Row 33
$ws.Range("F33").Value = "13-12-2022"
$ws.Range("G33").Formula = "'0"

# Row 34
$ws.Range("F34").Value = "13-12-2022"
$ws.Range("G34").Formula = "'0"

# Row 35
$ws.Range("F35").Value = "13-12-2022"
$ws.Range("G35").Formula = "'0"

# Row 36
$ws.Range("F36").Value = "13-12-2022"
$ws.Range("G36").Formula = "'0"

# Row 37
$ws.Range("F37").Value = "13-12-2022"
$ws.Range("G37").Formula = "'0"

# Row 38
$ws.Range("F38").Value = "13-12-2022"
$ws.Range("G38").Formula = "'0"

# Row 39
$ws.Range("F39").Value = "13-12-2022"
$ws.Range("G39").Formula = "'0"

# Row 40
$ws.Range("D40").Formula = "'0.04738"
$ws.Range("F40").Value = "13-12-2022"
$ws.Range("G40").Formula = "'0"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Formula = "'0.007084"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("F41").Value = "13-12-2022"
$ws.Range("G41").Formula = "'0"

# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Formula = "'0.004799"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("F42").Value = "13-12-2022"
$ws.Range("G42").Formula = "'0"

# Row 43
$ws.Range("D43").Formula = "'0.1119"
$ws.Range("F43").Value = "13-12-2022"
$ws.Range("G43").Formula = "'0"

# Row 44
$ws.Range("D44").Formula = "'0.01161"
$ws.Range("F44").Value = "13-12-2022"
$ws.Range("G44").Formula = "'0"

# Row 45
$ws.Range("D45").Formula = "'0.00006350"
$ws.Range("F45").Value = "13-12-2022"
$ws.Range("G45").Formula = "'0"

# Row 46
$ws.Range("D46").Formula = "'0.0009901"
$ws.Range("E46").Value = "45ACDXExchangeACXTBestin24h"
$ws.Range("F46").Value = "13-12-2022"
$ws.Range("G46").Formula = "'0"

# Row 47
$ws.Range("F47").Value = "13-12-2022"
$ws.Range("G47").Formula = "'0"

# Row 48
$ws.Range("D48").Formula = "'0.7231"
$ws.Range("F48").Value = "13-12-2022"
$ws.Range("G48").Formula = "'0"

# Row 49
$ws.Range("D49").Formula = "'0.001400"
$ws.Range("F49").Value = "13-12-2022"
$ws.Range("G49").Formula = "'0"

# Row 50
$ws.Range("D50").Formula = "'0.00001400"
$ws.Range("F50").Value = "13-12-2022"
$ws.Range("G50").Formula = "'0"

# Row 51
$ws.Range("F51").Value = "13-12-2022"
$ws.Range("G51").Formula = "'0"
